$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prediction ("Lookup") values for rows 28-85 (column C), per the commit diff.
$newC = @{
    28 = 0
    29 = 0
    30 = 0.011
    31 = 0.013
    32 = 0.016
    33 = 0.025
    34 = 0.039
    35 = 0.1
    36 = 0.126
    37 = 0.112
    38 = 0.134
    39 = 0.156
    40 = 0.164
    41 = 0.149
    42 = 0.284
    43 = 0.368
    44 = 0.478
    45 = 0.577
    46 = 0.625
    47 = 0.652
    48 = 0.67
    49 = 0.7
    50 = 0.708
    51 = 0.714
    52 = 0.712
    53 = 0.698
    54 = 0.688
    55 = 0.667
    56 = 0.647
    57 = 0.629
    58 = 0.63
    59 = 0.617
    60 = 0.594
    61 = 0.569
    62 = 0.539
    63 = 0.501
    64 = 0.487
    65 = 0.456
    66 = 0.419
    67 = 0.375
    68 = 0.342
    69 = 0.314
    70 = 0.29
    71 = 0.246
    72 = 0.198
    73 = 0.153
    74 = 0.126
    75 = 0.099
    76 = 0.076
    77 = 0.057
    78 = 0.042
    79 = 0.03
    80 = 0.022
    81 = 0.015
    82 = 0
    83 = 0
    84 = 0
    85 = 0
}

for ($r = 2; $r -le 96; $r++) {
    # Column A: the interval timestamp moves forward by exactly 30 days
    # (30.07.2024 -> 29.08.2024), keeping the same time-of-day fraction.
    $a = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $a + 30

    # Column D: the "Lookup" label is the new date prefix + the Interval number.
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 4).Value = "29.08.2024" + $b

    # Column C: updated prediction values where they changed.
    if ($newC.ContainsKey($r)) {
        $ws.Cells.Item($r, 3).Value = $newC[$r]
    }
}
